$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.229.07'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.38%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.860.00'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.07%  '

$ws.Range('E4').Value = '  +1.60%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.56'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.09%  '

$ws.Range('E6').Value = '  +1.46%  '

$ws.Range('E7').Value = '  +2.05%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3722'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.72%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07323'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.59%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9356'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.22%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.33'
$ws.Range('D11').ClearFormats()

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07865'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.39%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.853.11'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.97%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.421'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.63%  '

$ws.Range('E15').Value = '  +2.34%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '90.04'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.20%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.021'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.30%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008731'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.24%  '

$ws.Range('E19').Value = '  +1.41%  '

$ws.Range('E20').Value = '  +2.26%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '27.261.94'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.36%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.106'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.07%  '

$ws.Range('E23').Value = '  +0.74%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.954'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.37%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.79'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.38%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '18.50'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.65%  '

$ws.Range('E27').Value = '  -0.31%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '115.73'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.54%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.989'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.34%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.08887'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.87%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.348'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.37%  '

$ws.Range('E32').Value = '  +1.03%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.586'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.55%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7403'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.67%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.681'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.69%  '

$ws.Range('E36').Value = '  +3.60%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02031'
$ws.Range('D37').ClearFormats()

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05260'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.19%  '

$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5334'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.76%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.113'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.21%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1528'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.42%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.324'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.20%  '

$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '10.59'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.88%  '

$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4789'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.13%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.020'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.53%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '102.71'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.33%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.635'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.69%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '66.46'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.19%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06077'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.74%  '

$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.9015'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.38%  '

$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '36.69'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.41%  '
